# Apply the "updates aan the template" edit described by the commit diff.
$wb = $excel.ActiveWorkbook

$wsAlgemeen   = $wb.Worksheets.Item("Algemeen")
$wsGWE        = $wb.Worksheets.Item("GWE_Detail")
$wsSchoonmaak = $wb.Worksheets.Item("Schoonmaak")
$wsSchade     = $wb.Worksheets.Item("Schade")

# --- Algemeen ---------------------------------------------------------
# "Naam Klant *" test data swapped out for the new sample value.
$wsAlgemeen.Range("Klantnaam").Value = "underuse test"

# Schoonmaak pakket switches from "Intensief Schoonmaak" to "Basis Schoonmaak",
# which flips the "Inbegrepen Uren" formula result from 7 to 5.
$wsAlgemeen.Range("Schoonmaak_pakket").Value = "Basis Schoonmaak"

# --- GWE_Detail --------------------------------------------------------
# Updated meter end-readings (drives KWh/Gas verbruik + cost formulas).
$wsGWE.Range("KWh_eind").Value = 15020
$wsGWE.Range("Gas_eind").Value = 8050

# Row 19 was a blank example row; it now holds a real "water verbruik" line.
$wsGWE.Range("A19").Value = "super water verbruik"
$wsGWE.Range("B19").Value = 60
$wsGWE.Range("C19").Value = 4.5
$wsGWE.Range("D19").Formula = "=B19*C19"

# --- Schoonmaak ---------------------------------------------------------
# Updated hours actually worked.
$wsSchoonmaak.Range("Totaal_uren_gew").Value = 7

# --- Schade --------------------------------------------------------------
# Remove the bulk of the (example/demo) damage line items, keeping only the
# first two rows of real data; the shared formula in column D stays behind
# and simply evaluates to blank for the cleared rows.
$wsSchade.Range("A8:C22").ClearContents()
